$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4978.6
$ws.Range("I86").Value = 4973.5
$ws.Range("J86").Value = 4999
$ws.Range("K86").Value = 4973.5
$ws.Range("L86").Value = 4999
$ws.Range("M86").Value = -3850.5
$ws.Range("N86").Value = -7245
$ws.Range("H87").Value = 93720
$ws.Range("I87").Value = 90000
$ws.Range("J87").Value = 94960
$ws.Range("K87").Value = 90000
$ws.Range("L87").Value = 94960
$ws.Range("M87").Value = -88752
$ws.Range("N87").Value = -97456
$ws.Range("H89").Value = 4978.6
$ws.Range("I89").Value = 4973.5
$ws.Range("J89").Value = 4999
$ws.Range("K89").Value = 24867.5
$ws.Range("L89").Value = 24995
$ws.Range("M89").Value = -19251.5
$ws.Range("N89").Value = -36227
$ws.Range("H90").Value = 93720
$ws.Range("I90").Value = 90000
$ws.Range("J90").Value = 94960
$ws.Range("K90").Value = 270000
$ws.Range("L90").Value = 284880
$ws.Range("M90").Value = -263760
$ws.Range("N90").Value = -297360
$ws.Range("H121").Value = 1909.6
$ws.Range("J121").Value = 1909.6
$ws.Range("L121").Value = 5728.799999999999
$ws.Range("N121").Value = -9222.799999999999
$ws.Range("H125").Value = 4255
$ws.Range("I125").Value = 2765
$ws.Range("J125").Value = 5000
$ws.Range("K125").Value = 24885
$ws.Range("L125").Value = 45000
$ws.Range("M125").Value = -22425
$ws.Range("N125").Value = -49920
$ws.Range("H135").Value = 979.38464
$ws.Range("I135").Value = 979.38464
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8814.46176
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = $null
$ws.Range("N135").Value = -6279.46176
$ws.Range("H138").Value = 2684.8928
$ws.Range("I138").Value = 1121.069
$ws.Range("K138").Value = 3363.207
$ws.Range("M138").Value = 1776.793
$ws.Range("H141").Value = 2432.0435
$ws.Range("I141").Value = 1764.4
$ws.Range("J141").Value = 6883
$ws.Range("K141").Value = 5293.200000000001
$ws.Range("L141").Value = 20649
$ws.Range("M141").Value = -113.2000000000007
$ws.Range("N141").Value = -31009

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4306.5454
$ws.Range("I32").Value = 3972.7188
$ws.Range("K32").Value = 3972.7188
$ws.Range("M32").Value = -3685.7188
$ws.Range("H61").Value = 2722.0667
$ws.Range("I61").Value = 2702.4285
$ws.Range("K61").Value = 2702.4285
$ws.Range("M61").Value = -2490.4285
$ws.Range("H102").Value = 2951.75
$ws.Range("I102").Value = 2302
$ws.Range("K102").Value = 2302
$ws.Range("M102").Value = -680
$ws.Range("H132").Value = 2037.0233
$ws.Range("I132").Value = 1883.1538
$ws.Range("K132").Value = 5649.4614
$ws.Range("M132").Value = -3119.4614
$ws.Range("H136").Value = 2722.0667
$ws.Range("I136").Value = 2702.4285
$ws.Range("K136").Value = 8107.2855
$ws.Range("M136").Value = -5557.2855

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 35000
$ws.Range("J35").Value = 35000
$ws.Range("L35").Value = 35000
$ws.Range("N35").Value = -35620
$ws.Range("H86").Value = 2849
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = $null
$ws.Range("H89").Value = 2849
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = $null
$ws.Range("H134").Value = 2859
$ws.Range("I134").Value = 2859
$ws.Range("K134").Value = 8577
$ws.Range("M134").Value = -6042

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3249
$ws.Range("I58").Value = 3299
$ws.Range("K58").Value = 3299
$ws.Range("M58").Value = -3096
$ws.Range("H86").Value = 19349.777
$ws.Range("I86").Value = 5197.3
$ws.Range("K86").Value = 5197.3
$ws.Range("M86").Value = -4074.3
$ws.Range("H89").Value = 19349.777
$ws.Range("I89").Value = 5197.3
$ws.Range("K89").Value = 25986.5
$ws.Range("M89").Value = -20370.5
$ws.Range("H132").Value = 3820.9333
$ws.Range("I132").Value = 3581.9
$ws.Range("J132").Value = 4299
$ws.Range("K132").Value = 10745.7
$ws.Range("L132").Value = 12897
$ws.Range("M132").Value = -8215.700000000001
$ws.Range("N132").Value = -17957
$ws.Range("H134").Value = 5979.2856
$ws.Range("I134").Value = 5979.2856
$ws.Range("K134").Value = 17937.8568
$ws.Range("M134").Value = -15402.8568
$ws.Range("H136").Value = 3249
$ws.Range("I136").Value = 3299
$ws.Range("K136").Value = 9897
$ws.Range("M136").Value = -7347

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 9999
$ws.Range("J29").Value = 9999
$ws.Range("L29").Value = 29997
$ws.Range("N29").Value = -30551
$ws.Range("H68").Value = 1600
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = $null
$ws.Range("H71").Value = 1600
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = $null
$ws.Range("H92").Value = 1500
$ws.Range("J92").Value = 2000
$ws.Range("L92").Value = 6000
$ws.Range("N92").Value = -8496
$ws.Range("H122").Value = 1225.75
$ws.Range("I122").Value = 503
$ws.Range("J122").Value = 1466.6666
$ws.Range("K122").Value = 4527
$ws.Range("L122").Value = 13199.9994
$ws.Range("M122").Value = -2077
$ws.Range("N122").Value = -18099.9994
$ws.Range("H131").Value = 1383.5625
$ws.Range("J131").Value = 1811.625
$ws.Range("L131").Value = 5434.875
$ws.Range("N131").Value = -15514.875
$ws.Range("H132").Value = 4129.25
$ws.Range("I132").Value = 3607.6
$ws.Range("K132").Value = 32468.4
$ws.Range("M132").Value = -29938.4
$ws.Range("H140").Value = 3322.3333
$ws.Range("I140").Value = 3322.3333
$ws.Range("K140").Value = 9966.999899999999
$ws.Range("M140").Value = -4786.999899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4476.7
$ws.Range("I80").Value = 3862.3333
$ws.Range("J80").Value = 10006
$ws.Range("K80").Value = 3862.3333
$ws.Range("L80").Value = 10006
$ws.Range("M80").Value = -2864.3333
$ws.Range("N80").Value = -12002
$ws.Range("H83").Value = 4476.7
$ws.Range("I83").Value = 3862.3333
$ws.Range("J83").Value = 10006
$ws.Range("K83").Value = 19311.6665
$ws.Range("L83").Value = 50030
$ws.Range("M83").Value = -14319.6665
$ws.Range("N83").Value = -60014

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 209.2
$ws.Range("I55").Value = 99.25
$ws.Range("J55").Value = 282.5
$ws.Range("K55").Value = 99.25
$ws.Range("L55").Value = 282.5
$ws.Range("M55").Value = 73.75
$ws.Range("N55").Value = -628.5
$ws.Range("H93").Value = 2119.3333
$ws.Range("I93").Value = 2071.875
$ws.Range("K93").Value = 2071.875
$ws.Range("M93").Value = -823.875
$ws.Range("H122").Value = 2779.3572
$ws.Range("I122").Value = 2608.7693
$ws.Range("K122").Value = 7826.3079
$ws.Range("M122").Value = -5376.3079
$ws.Range("H132").Value = 2587.4443
$ws.Range("I132").Value = 2059.6
$ws.Range("J132").Value = 3247.25
$ws.Range("K132").Value = 6178.799999999999
$ws.Range("L132").Value = 9741.75
$ws.Range("M132").Value = -3648.799999999999
$ws.Range("N132").Value = -14801.75
$ws.Range("H136").Value = 3202
$ws.Range("I136").Value = 3045.0715
$ws.Range("K136").Value = 9135.2145
$ws.Range("M136").Value = -6585.2145

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 28163
$ws.Range("J45").Value = 28163
$ws.Range("L45").Value = 28163
$ws.Range("N45").Value = -29145
$ws.Range("H126").Value = 2536.4546
$ws.Range("I126").Value = 2576.25
$ws.Range("K126").Value = 7728.75
$ws.Range("M126").Value = -5258.75
$ws.Range("H132").Value = 1791.6842
$ws.Range("I132").Value = 1626.6
$ws.Range("K132").Value = 4879.799999999999
$ws.Range("M132").Value = -2349.799999999999
$ws.Range("H136").Value = 1163.32
$ws.Range("I136").Value = 948.0454999999999
$ws.Range("J136").Value = 2742
$ws.Range("K136").Value = 2844.1365
$ws.Range("L136").Value = 8226
$ws.Range("M136").Value = -294.1364999999996
$ws.Range("N136").Value = -13326
